# Update so MASTER is read in directly.
# The MASTER-derived rows for "Entiat River Lake 06" and "Entiat River Lake 07"
# (the original data rows 9 and 10) are no longer present in the refreshed
# MASTER extract, so they are removed and every following row shifts up by two,
# shrinking the used range from A1:W31 down to A1:W29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 and 10 entirely; Excel shifts the remaining rows (old 11-31) up.
$ws.Rows("9:10").Delete()
